$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04976167759201644
$ws.Range("D2").Value = 0.02354875200813211
$ws.Range("E2").Value = 0.1064597284362492
$ws.Range("F2").Value = 0.4444044161742013
$ws.Range("G2").Value = 0.2896642802316904
$ws.Range("H2").Value = 0.4621988561487242
$ws.Range("K2").Value = 1.185113622351707
$ws.Range("M2").Value = 0.3745258566448868
$ws.Range("O2").Value = 1.429091374701287

$ws.Range("C3").Value = 0.04416796813832491
$ws.Range("D3").Value = 0.02094628914163366
$ws.Range("E3").Value = 0.1030373694395337
$ws.Range("F3").Value = 0.4453197449422888
$ws.Range("G3").Value = 0.2917212057711183
$ws.Range("H3").Value = 0.4677580831142123
$ws.Range("K3").Value = 1.034403058535588
$ws.Range("M3").Value = 0.3320957871796466
$ws.Range("O3").Value = 1.444796239223365

$ws.Range("C4").Value = 0.04074669862028202
$ws.Range("D4").Value = 0.01934040412479732
$ws.Range("E4").Value = 0.1010777977424908
$ws.Range("F4").Value = 0.446329959083954
$ws.Range("G4").Value = 0.2933611242527974
$ws.Range("H4").Value = 0.471497971049935
$ws.Range("K4").Value = 0.9414913515508374
$ws.Range("M4").Value = 0.3060529549820217
$ws.Range("O4").Value = 1.455911708282727

$ws.Range("C5").Value = 0.03935584826531624
$ws.Range("D5").Value = 0.01868403836257215
$ws.Range("E5").Value = 0.1003145679918411
$ws.Range("F5").Value = 0.4468540151482401
$ws.Range("G5").Value = 0.2941237627368096
$ws.Range("H5").Value = 0.4731039981795462
$ws.Range("K5").Value = 0.9035373628456966
$ws.Range("M5").Value = 0.2954429059795203
$ws.Range("O5").Value = 1.460810506897573

$ws.Range("C6").Value = 0.03912510072476039
$ws.Range("D6").Value = 0.01857493268487787
$ws.Range("E6").Value = 0.1001899571268901
$ws.Range("F6").Value = 0.4469478127425859
$ws.Range("G6").Value = 0.2942560846409634
$ws.Range("H6").Value = 0.4733756267139455
$ws.Range("K6").Value = 0.897229667827105
$ws.Range("M6").Value = 0.2936812810339191
$ws.Range("O6").Value = 1.461646211027173

$ws.Range("C7").Value = 0.04072792757982313
$ws.Range("D7").Value = 0.01933155999722658
$ws.Range("E7").Value = 0.1010673620303848
$ws.Range("F7").Value = 0.4463365720650572
$ws.Range("G7").Value = 0.2933710279812658
$ws.Range("H7").Value = 0.4715192986547478
$ws.Range("K7").Value = 0.9409798587489888
$ws.Range("M7").Value = 0.3059098532161642
$ws.Range("O7").Value = 1.455976282027933

$ws.Range("C8").Value = 0.04783020800803683
$ws.Range("D8").Value = 0.02265310304994017
$ws.Range("E8").Value = 0.1052500700674486
$ws.Range("F8").Value = 0.4446268059206275
$ws.Range("G8").Value = 0.2902950113549565
$ws.Range("H8").Value = 0.4640478473742178
$ws.Range("K8").Value = 1.133227879847993
$ws.Range("M8").Value = 0.3598941541545457
$ws.Range("O8").Value = 1.434200080156117

$ws.Range("C9").Value = 0.06186381656110029
$ws.Range("D9").Value = 0.02910168481030695
$ws.Range("E9").Value = 0.1145929958443261
$ws.Range("F9").Value = 0.4448451754501406
$ws.Range("G9").Value = 0.2872736027031451
$ws.Range("H9").Value = 0.4519918980135174
$ws.Range("K9").Value = 1.507159020275026
$ws.Range("M9").Value = 0.4658262518548639
$ws.Range("O9").Value = 1.40323384413459

$ws.Range("C10").Value = 0.07224084153885713
$ws.Range("D10").Value = 0.03379793476456427
$ws.Range("E10").Value = 0.1221751907550441
$ws.Range("F10").Value = 0.4472043449687817
$ws.Range("G10").Value = 0.2869166850592819
$ws.Range("H10").Value = 0.4447235213204692
$ws.Range("K10").Value = 1.779919698878075
$ws.Range("M10").Value = 0.543697373624866
$ws.Range("O10").Value = 1.387711108322208

$ws.Range("C11").Value = 0.07697652628792184
$ws.Range("D11").Value = 0.03592499005702621
$ws.Range("E11").Value = 0.1257853289219355
$ws.Range("F11").Value = 0.4487597052179098
$ws.Range("G11").Value = 0.2871647399685031
$ws.Range("H11").Value = 0.441763484804369
$ws.Range("K11").Value = 1.903560113592903
$ws.Range("M11").Value = 0.5791329993371335
$ws.Range("O11").Value = 1.382234645339508

$ws.Range("C12").Value = 0.07877199435897353
$ws.Range("D12").Value = 0.03672907349348975
$ws.Range("E12").Value = 0.1271759211881047
$ws.Range("F12").Value = 0.4494183707301289
$ws.Range("G12").Value = 0.2873181578497679
$ws.Range("H12").Value = 0.4406925285644121
$ws.Range("K12").Value = 1.950314179475697
$ws.Range("M12").Value = 0.5925531113533964
$ws.Range("O12").Value = 1.380390015799264

$ws.Range("C13").Value = 0.07838521193345116
$ws.Range("D13").Value = 0.03655596218318635
$ws.Range("E13").Value = 0.1268753807651279
$ws.Range("F13").Value = 0.4492734105465246
$ws.Range("G13").Value = 0.2872824635405777
$ws.Range("H13").Value = 0.4409209546625021
$ws.Range("K13").Value = 1.940247826331984
$ws.Range("M13").Value = 0.5896627916236667
$ws.Range("O13").Value = 1.380777076234864

$ws.Range("C14").Value = 0.07712419714138719
$ws.Range("D14").Value = 0.03599117055264855
$ws.Range("E14").Value = 0.1258992604874436
$ws.Range("F14").Value = 0.4488124954508024
$ws.Range("G14").Value = 0.2871761669156001
$ws.Range("H14").Value = 0.4416743748348466
$ws.Range("K14").Value = 1.907407933968102
$ws.Range("M14").Value = 0.580237053291782
$ws.Range("O14").Value = 1.382078285969044

$ws.Range("C15").Value = 0.07635207114699938
$ws.Range("D15").Value = 0.03564503736497215
$ws.Range("E15").Value = 0.1253044313550404
$ws.Range("F15").Value = 0.4485392569554136
$ws.Range("G15").Value = 0.2871188176809767
$ws.Range("H15").Value = 0.4421423748251101
$ws.Range("K15").Value = 1.887283894992947
$ws.Range("M15").Value = 0.5744636921273099
$ws.Range("O15").Value = 1.382905199479069

$ws.Range("C16").Value = 0.07193165707916194
$ws.Range("D16").Value = 0.03365873542630027
$ws.Range("E16").Value = 0.1219425298844428
$ws.Range("F16").Value = 0.4471124308232532
$ws.Range("G16").Value = 0.2869087727400768
$ws.Range("H16").Value = 0.4449239484150382
$ws.Range("K16").Value = 1.771830431433614
$ws.Range("M16").Value = 0.5413817838433914
$ws.Range("O16").Value = 1.388101010010359

$ws.Range("C17").Value = 0.06922374442109458
$ws.Range("D17").Value = 0.03243778779874162
$ws.Range("E17").Value = 0.1199215977134998
$ws.Range("F17").Value = 0.44636086223219
$ws.Range("G17").Value = 0.2868853680586838
$ws.Range("H17").Value = 0.4467191677336757
$ws.Range("K17").Value = 1.700888946782072
$ws.Range("M17").Value = 0.521089899662357
$ws.Range("O17").Value = 1.391695316886683

$ws.Range("C18").Value = 0.06766764871424868
$ws.Range("D18").Value = 0.03173465933524255
$ws.Range("E18").Value = 0.1187743393839469
$ws.Range("F18").Value = 0.445973942456483
$ws.Range("G18").Value = 0.2869105090387478
$ws.Range("H18").Value = 0.4477843256563006
$ws.Range("K18").Value = 1.660044041992023
$ws.Range("M18").Value = 0.5094196579132273
$ws.Range("O18").Value = 1.393911769883431

$ws.Range("C19").Value = 0.06714102625613805
$ws.Range("D19").Value = 0.03149644393731421
$ws.Range("E19").Value = 0.1183884846170287
$ws.Range("F19").Value = 0.4458507184405534
$ws.Range("G19").Value = 0.2869256371322351
$ws.Range("H19").Value = 0.4481505639561973
$ws.Range("K19").Value = 1.646207656965203
$ws.Range("M19").Value = 0.5054685140228088
$ws.Range("O19").Value = 1.394687792237065

$ws.Range("C20").Value = 0.06951185901989732
$ws.Range("D20").Value = 0.03256785033791232
$ws.Range("E20").Value = 0.1201351610155754
$ws.Range("F20").Value = 0.4464361705900046
$ws.Range("G20").Value = 0.2868838608979303
$ws.Range("H20").Value = 0.4465246891008192
$ws.Range("K20").Value = 1.708445076212797
$ws.Range("M20").Value = 0.5232498921043458
$ws.Range("O20").Value = 1.391297255914225

$ws.Range("C21").Value = 0.07749452900496578
$ws.Range("D21").Value = 0.03615710156400098
$ws.Range("E21").Value = 0.1261853296689068
$ws.Range("F21").Value = 0.4489459833954044
$ws.Range("G21").Value = 0.2872057706956639
$ws.Range("H21").Value = 0.4414517204846362
$ws.Range("K21").Value = 1.917055615108382
$ws.Range("M21").Value = 0.5830055849719571
$ws.Range("O21").Value = 1.381689858848347

$ws.Range("C22").Value = 0.08272429259136516
$ws.Range("D22").Value = 0.03849478176951493
$ws.Range("E22").Value = 0.1302766765188679
$ws.Range("F22").Value = 0.4509926257408949
$ws.Range("G22").Value = 0.287763096096441
$ws.Range("H22").Value = 0.4384274246993556
$ws.Range("K22").Value = 2.053009485300038
$ws.Range("M22").Value = 0.6220674804780089
$ws.Range("O22").Value = 1.376747344056298

$ws.Range("C23").Value = 0.07993191659780052
$ws.Range("D23").Value = 0.03724787567604437
$ws.Range("E23").Value = 0.1280803718771892
$ws.Range("F23").Value = 0.4498629994698575
$ws.Range("G23").Value = 0.2874337422799158
$ws.Range("H23").Value = 0.4400148603112939
$ws.Range("K23").Value = 1.980484452676876
$ws.Range("M23").Value = 0.6012187579339496
$ws.Range("O23").Value = 1.379262547931802

$ws.Range("C24").Value = 0.06938160016970585
$ws.Range("D24").Value = 0.03250905276507865
$ws.Range("E24").Value = 0.1200385635906542
$ws.Range("F24").Value = 0.4464019830151997
$ws.Range("G24").Value = 0.2868844220943245
$ws.Range("H24").Value = 0.4466125099116596
$ws.Range("K24").Value = 1.705029135862731
$ws.Range("M24").Value = 0.5222733724871631
$ws.Range("O24").Value = 1.391476752013034

$ws.Range("C25").Value = 0.05805580429729673
$ws.Range("D25").Value = 0.02736432717053106
$ws.Range("E25").Value = 0.1119410409154469
$ws.Range("F25").Value = 0.4444015346785477
$ws.Range("G25").Value = 0.2877658093292581
$ws.Range("H25").Value = 0.4549747754221656
$ws.Range("K25").Value = 1.406339192323969
$ws.Range("M25").Value = 0.4371613667730117
$ws.Range("O25").Value = 1.41034695273207
